# Update symbol list data on the active sheet.
#
# Column D holds numeric-looking price strings that must stay stored as
# TEXT (the source workbook uses inlineStr for every data cell, not
# numbers). Assigning a numeric-looking string straight to .Value lets
# Excel auto-convert it to a real number, so for those cells we force
# text formatting first, then clear the formatting delta back out again
# so no stray NumberFormat/style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# --- Simple price (column D) updates -----------------------------------
Set-TextValue "D2"  "244.59"
Set-TextValue "D3"  "21.87"
Set-TextValue "D4"  "5.387"
Set-TextValue "D5"  "0.06006"
Set-TextValue "D6"  "3.391"
Set-TextValue "D7"  "0.8169"
Set-TextValue "D8"  "0.9520"
Set-TextValue "D11" "0.03283"
Set-TextValue "D12" "0.03052"
Set-TextValue "D16" "0.04809"
Set-TextValue "D17" "0.0005912"
Set-TextValue "D18" "0.005506"
Set-TextValue "D19" "0.004158"
Set-TextValue "D20" "0.0009866"
Set-TextValue "D22" "6.416"
Set-TextValue "D26" "0.00007002"
Set-TextValue "D40" "0.03997"
Set-TextValue "D44" "0.005817"
Set-TextValue "D45" "0.00005133"
Set-TextValue "D47" "0.8603"
Set-TextValue "D48" "0.004079"

# --- Row 13/14 swap: BitMartToken <-> MCDex -----------------------------
$ws.Range("B13").Value = "MCDex"
$ws.Range("C13").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D13" "4.010"
$ws.Range("E13").Value = "12MCDexMCB"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09409"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# --- Row 41/42/43 rotation: KickToken -> BKEXToken -> CEJI -> KickToken -
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1075"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002681"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003041"
$ws.Range("E43").Value = "42KickTokenKICK"
